$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet had an editing mistake: an orphan row (old row 13)
# holding only "2143261 - André Luis Ferraz" in columns B/C, sitting under
# the "Docentes responsáveis:" label instead of being paired with it, plus
# several description cells further down were paired with the wrong label.
# Fix: delete the stray row (shifting everything below it up by one, which
# also restores the correct per-row heights further down to match), then
# correct the handful of description cells that still hold mismatched text.

# Remove the stray row 13 - shifts rows 14+ up to 13+, row heights follow.
$ws.Rows.Item(13).Delete()

# After the shift, patch up the description cells that now hold the wrong
# text so each column-A label is paired with its correct B/C content.
$ws.Range("B10").Value = "2143261 - André Luis Ferraz"
$ws.Range("C10").Value = "2143261 - André Luis Ferraz"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2018" must stay literal text (not get auto-parsed into a date
# serial), so clone it from the already-correct text cell via copy/paste
# of values only (also keeps the style table untouched).
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

$ws.Range("B18").Value = "2143261 - André Luis Ferraz"
$ws.Range("C18").Value = "2143261 - André Luis Ferraz"

$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."

$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
